$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New secretChildName / secretChildEmailID pairings for rows 2-14
$values = @(
    @("Spencer Allen", "spencer.allen@acme.com"),
    @("Hamish Murray", "hamish.murray@acme.com"),
    @("Charlie Ross", "charlie.ross.jr@acme.com"),
    @("Mark Lawrence", "mark.lawrence@acme.com"),
    @("Layla Graham", "layla.graham@acme.com"),
    @("Isabella Scott", "isabella.scott@acme.com"),
    @("Matthew King", "matthew.king@acme.com"),
    @("Charlie Ross", "charlie.ross@acme.com"),
    @("Piper Stewart", "piper.stewart@acme.com"),
    @("Matthew King", "matthew.king.jr@acme.com"),
    @("Charlie Wright", "charlie.wright@acme.com"),
    @("Ethan Murray", "ethan.murray@acme.com"),
    @("Benjamin Collins", "benjamin.collins@acme.com")
)

$row = 2
foreach ($pair in $values) {
    $ws.Cells.Item($row, 3).Value = $pair[0]
    $ws.Cells.Item($row, 4).Value = $pair[1]
    $row++
}
